# Ders 0.pptx - slide 7 ("Grafiksel Kullanıcı Arayüzü"):
# The content placeholder's first paragraph ends with a single run:
#   " - GIU): Web tarayıcısında gördüğümüz grafikler."
# It must become three runs (fixing the "GIU" -> "GUI" typo along the way):
#   " "         (unchanged formatting)
#   "- GUI): "  (the corrected acronym)
#   "Web tarayıcısında gördüğümüz grafikler."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$full = $tr.Text
$idx = $full.IndexOf(" - GIU): Web tarayıcısında gördüğümüz grafikler.")

# 1-based start of the run inside the overall shape TextRange.
$start = $idx + 1

# Piece 1: the leading space.
$part1 = $tr.Characters($start, 1)
$part1.Text = " "

# Piece 2: "- GIU): " -> "- GUI): " (typo fix).
$part2 = $tr.Characters($start + 1, 8)
$part2.Text = "- GUI): "

# Piece 3: the remaining, unchanged sentence.
$part3 = $tr.Characters($start + 9, 39)
$part3.Text = "Web tarayıcısında gördüğümüz grafikler."
